$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest cryptos snapshot.
# Numeric-looking price strings (single "." as a plausible decimal point) are
# forced through a Text number-format round-trip so Excel keeps them as the
# literal text the source feed provides (e.g. "241.96", "8.080") instead of
# silently parsing them into actual numbers; the temporary formatting is then
# cleared so the cell is left with no explicit style, matching the original.

$ws.Range("D2").Value = "29.513.55"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").Value = "1.878.22"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("E5").Value = "  +2.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.56%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07901"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3098"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.43"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08276"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.34%  "

$ws.Range("D12").Value = "1.892.03"
$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7293"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.281"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.25"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.89%  "

$ws.Range("D16").Value = "29.520.09"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.904"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "246.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007862"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.33"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("D21").Value = "2.126.76"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.080"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.82%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("E25").Value = "  +14.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.055"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.42%  "

$ws.Range("E29").Value = "  -3.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.495"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.398"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.110"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05221"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.951"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.199"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7286"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.681"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("E38").Value = "  +1.19%  "

$ws.Range("D39").Value = "1.203.61"
$ws.Range("E39").Value = "  +5.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.707"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9111"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.182"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.36%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.44"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").Value = "2.022.61"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5297"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("E48").Value = "  +3.17%  "

$ws.Range("E49").Value = "  +9.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.320"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4328"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.87%  "
